# Add review UI & AddReviewController
# The "Movie Review List" row (row 10) gains its first review entry,
# encoded as "|rating:<n>review:<text>" in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "|rating:5review:Great!"
